$d = $word.ActiveDocument

$pairs = @(
    @("50×34=1700", "22×73=1606"),
    @("47×93=4371", "14×19=266"),
    @("52×23=1196", "13×49=637"),
    @("47×26=1222", "93×48=4464"),
    @("63×82=5166", "23×22=506"),
    @("24×16=384", "59×92=5428"),
    @("94×30=2820", "53×88=4664"),
    @("94×12=1128", "13×79=1027"),
    @("74×58=4292", "63×98=6174"),
    @("28×50=1400", "13×39=507"),
    @("80×37=2960", "91×29=2639"),
    @("63×86=5418", "12×91=1092"),
    @("54×21=1134", "15×35=525"),
    @("53×40=2120", "60×69=4140"),
    @("54×80=4320", "25×22=550"),
    @("90×88=7920", "74×11=814"),
    @("14×34=476", "88×30=2640"),
    @("12×86=1032", "25×17=425"),
    @("18×72=1296", "49×67=3283"),
    @("21×63=1323", "34×31=1054"),
    @("43×45=1935", "88×69=6072"),
    @("47×34=1598", "81×78=6318"),
    @("97×37=3589", "33×45=1485"),
    @("31×17=527", "89×90=8010"),
    @("33×21=693", "16×48=768")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
